$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.843543291091919
$ws.Range("B1").Value = 2.554818630218506
$ws.Range("C1").Value = 2.785227298736572
$ws.Range("D1").Value = 3.336065292358398
$ws.Range("E1").Value = 0.9857743978500366
